$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 579
$ws.Range("J31").Value = 579
$ws.Range("L31").Value = 1737
$ws.Range("N31").Value = -2197
$ws.Range("H39").Value = 656.4286
$ws.Range("I39").Value = 141.66667
$ws.Range("K39").Value = 425.00001
$ws.Range("M39").Value = -129.00001
$ws.Range("H47").Value = 10000
$ws.Range("I47").Value = 10000
$ws.Range("K47").Value = 10000
$ws.Range("M47").Value = -9028
$ws.Range("H70").Value = 2677.4443
$ws.Range("I70").Value = 900
$ws.Range("J70").Value = 3185.2856
$ws.Range("K70").Value = 2700
$ws.Range("L70").Value = 9555.856800000001
$ws.Range("M70").Value = -2430
$ws.Range("N70").Value = -10095.8568
$ws.Range("H73").Value = 2677.4443
$ws.Range("I73").Value = 900
$ws.Range("J73").Value = 3185.2856
$ws.Range("K73").Value = 2700
$ws.Range("L73").Value = 9555.856800000001
$ws.Range("M73").Value = -1764
$ws.Range("N73").Value = -11427.8568
$ws.Range("H80").Value = 299.08334
$ws.Range("I80").Value = 299.08334
$ws.Range("K80").Value = 897.2500200000001
$ws.Range("M80").Value = 100.7499799999999
$ws.Range("H81").Value = 95000
$ws.Range("J81").Value = 95000
$ws.Range("L81").Value = 95000
$ws.Range("N81").Value = -96996
$ws.Range("H83").Value = 299.08334
$ws.Range("I83").Value = 299.08334
$ws.Range("K83").Value = 2691.75006
$ws.Range("M83").Value = 2300.24994
$ws.Range("H84").Value = 95000
$ws.Range("J84").Value = 95000
$ws.Range("L84").Value = 285000
$ws.Range("N84").Value = -294984
$ws.Range("H92").Value = 592.8
$ws.Range("I92").Value = 541
$ws.Range("K92").Value = 541
$ws.Range("M92").Value = 707
$ws.Range("H107").Value = 1268.1818
$ws.Range("I107").Value = 1268.1818
$ws.Range("K107").Value = 1268.1818
$ws.Range("M107").Value = 651.8181999999999
$ws.Range("H112").Value = 3369.647
$ws.Range("J112").Value = 3369.647
$ws.Range("L112").Value = 10108.941
$ws.Range("N112").Value = -12324.941
$ws.Range("H113").Value = 4941.5557
$ws.Range("I113").Value = 2296.4285
$ws.Range("K113").Value = 2296.4285
$ws.Range("M113").Value = 957.5715
$ws.Range("H116").Value = 3741.8333
$ws.Range("I116").Value = 3451
$ws.Range("K116").Value = 3451
$ws.Range("M116").Value = -9
$ws.Range("H131").Value = 104024.3
$ws.Range("I131").Value = 104024.3
$ws.Range("K131").Value = 312072.9
$ws.Range("M131").Value = -307032.9
$ws.Range("H132").Value = 1974.6923
$ws.Range("I132").Value = 1701.7916
$ws.Range("J132").Value = 5249.5
$ws.Range("K132").Value = 5105.3748
$ws.Range("L132").Value = 15748.5
$ws.Range("M132").Value = -2575.3748
$ws.Range("N132").Value = -20808.5
$ws.Range("H137").Value = 1973.2858
$ws.Range("I137").Value = 1490.75
$ws.Range("K137").Value = 4472.25
$ws.Range("M137").Value = -1922.25
$ws.Range("H141").Value = 3159.7856
$ws.Range("I141").Value = 3211.8333
$ws.Range("K141").Value = 9635.499899999999
$ws.Range("M141").Value = -4455.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6916.3335
$ws.Range("I2").Value = 6700
$ws.Range("K2").Value = 6700
$ws.Range("M2").Value = -6587
$ws.Range("H32").Value = 5586.345
$ws.Range("I32").Value = 5419.7925
$ws.Range("K32").Value = 5419.7925
$ws.Range("M32").Value = -5132.7925
$ws.Range("H45").Value = 257289
$ws.Range("I45").Value = 405662.4
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 405662.4
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -405285.4
$ws.Range("N45").Value = -10754
$ws.Range("H102").Value = 3110.6875
$ws.Range("I102").Value = 2234.75
$ws.Range("K102").Value = 2234.75
$ws.Range("M102").Value = -612.75
$ws.Range("H110").Value = 5028.4585
$ws.Range("I110").Value = 3616.3125
$ws.Range("J110").Value = 7852.75
$ws.Range("K110").Value = 3616.3125
$ws.Range("L110").Value = 7852.75
$ws.Range("M110").Value = -1571.3125
$ws.Range("N110").Value = -11942.75
$ws.Range("H116").Value = 6916.3335
$ws.Range("I116").Value = 6700
$ws.Range("K116").Value = 6700
$ws.Range("M116").Value = -4406
$ws.Range("H132").Value = 5411.7856
$ws.Range("I132").Value = 4449.3335
$ws.Range("K132").Value = 13348.0005
$ws.Range("M132").Value = -10818.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6916.3335
$ws.Range("I3").Value = 6700
$ws.Range("K3").Value = 6700
$ws.Range("M3").Value = -6586
$ws.Range("H20").Value = 4720.7144
$ws.Range("J20").Value = 5499
$ws.Range("L20").Value = 5499
$ws.Range("N20").Value = -5993
$ws.Range("H75").Value = 9100
$ws.Range("I75").Value = 5466.6665
$ws.Range("K75").Value = 5466.6665
$ws.Range("M75").Value = -4530.6665
$ws.Range("H78").Value = 9100
$ws.Range("I78").Value = 5466.6665
$ws.Range("K78").Value = 16399.9995
$ws.Range("M78").Value = -11719.9995
$ws.Range("H100").Value = 25279.666
$ws.Range("J100").Value = 25279.666
$ws.Range("L100").Value = 25279.666
$ws.Range("N100").Value = -27443.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 519
$ws.Range("I107").Value = 446.7857
$ws.Range("K107").Value = 446.7857
$ws.Range("M107").Value = 1473.2143
$ws.Range("H122").Value = 5777.6665
$ws.Range("I122").Value = 5416.5
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 16249.5
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -13799.5
$ws.Range("N122").Value = -24400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 9258.546
$ws.Range("I2").Value = 82
$ws.Range("J2").Value = 12699.75
$ws.Range("K2").Value = 492
$ws.Range("L2").Value = 76198.5
$ws.Range("M2").Value = -379
$ws.Range("N2").Value = -76424.5
$ws.Range("H5").Value = 704
$ws.Range("I5").Value = 704
$ws.Range("K5").Value = 2112
$ws.Range("M5").Value = -2000
$ws.Range("H33").Value = 458.5
$ws.Range("J33").Value = 873.25
$ws.Range("L33").Value = 5239.5
$ws.Range("N33").Value = -5805.5
$ws.Range("H80").Value = 5999
$ws.Range("J80").Value = 5998.5
$ws.Range("L80").Value = 17995.5
$ws.Range("N80").Value = -19867.5
$ws.Range("H83").Value = 5999
$ws.Range("J83").Value = 5998.5
$ws.Range("L83").Value = 53986.5
$ws.Range("N83").Value = -63346.5
$ws.Range("H98").Value = 241.44444
$ws.Range("I98").Value = 240
$ws.Range("J98").Value = 242.16667
$ws.Range("K98").Value = 720
$ws.Range("L98").Value = 726.50001
$ws.Range("M98").Value = 778
$ws.Range("N98").Value = -3722.50001
$ws.Range("H128").Value = 404741.75
$ws.Range("I128").Value = 404741.75
$ws.Range("K128").Value = 1214225.25
$ws.Range("M128").Value = -1209245.25
$ws.Range("H135").Value = 704
$ws.Range("I135").Value = 704
$ws.Range("K135").Value = 6336
$ws.Range("M135").Value = -3801

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 30000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 30000
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 30000
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -31262
$ws.Range("H80").Value = 5001
$ws.Range("I80").Value = 5001
$ws.Range("K80").Value = 5001
$ws.Range("M80").Value = -4003
$ws.Range("H83").Value = 5001
$ws.Range("I83").Value = 5001
$ws.Range("K83").Value = 25005
$ws.Range("M83").Value = -20013
$ws.Range("H102").Value = 4153.923
$ws.Range("I102").Value = 3743.2
$ws.Range("J102").Value = 5523
$ws.Range("K102").Value = 3743.2
$ws.Range("L102").Value = 5523
$ws.Range("M102").Value = -2121.2
$ws.Range("N102").Value = -8767
$ws.Range("H122").Value = 3061.6924
$ws.Range("I122").Value = 2478.111
$ws.Range("K122").Value = 7434.333
$ws.Range("M122").Value = -4984.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3315
$ws.Range("I122").Value = 3265.3333
$ws.Range("K122").Value = 9795.999899999999
$ws.Range("M122").Value = -7345.999899999999
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 13269.4
$ws.Range("I132").Value = 14212.3125
$ws.Range("K132").Value = 42636.9375
$ws.Range("M132").Value = -40106.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 16672333
$ws.Range("I4").Value = 16672333
$ws.Range("K4").Value = 16672333
$ws.Range("M4").Value = -16672220
$ws.Range("H17").Value = 22888.334
$ws.Range("I17").Value = 22888.334
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 22888.334
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -22716.334
$ws.Range("N17").ClearContents()
$ws.Range("H122").Value = 3474.5186
$ws.Range("I122").Value = 2671.4119
$ws.Range("K122").Value = 8014.2357
$ws.Range("M122").Value = -5564.2357
$ws.Range("H132").Value = 2681.7097
$ws.Range("I132").Value = 2820
$ws.Range("K132").Value = 8460
$ws.Range("M132").Value = -5930
